$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.882.06"
$ws.Range("E2").Value = "  +9.50%  "
$ws.Range("D3").Value = "3.521.60"
$ws.Range("E3").Value = "  +10.93%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "191.91"
$ws.Range("E5").Value = "  +13.33%  "
$ws.Range("D6").Value = "557.29"
$ws.Range("E6").Value = "  +10.16%  "
$ws.Range("D7").Value = "3.515.58"
$ws.Range("E7").Value = "  +10.89%  "
$ws.Range("E8").Value = "  +5.19%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("E10").Value = "  +9.28%  "
$ws.Range("D11").Value = "57.28"
$ws.Range("E12").Value = "  +17.40%  "
$ws.Range("D13").Value = "0.0000278"
$ws.Range("E13").Value = "  +11.87%  "
$ws.Range("E14").Value = "  +8.52%  "
$ws.Range("D15").Value = "4.085.32"
$ws.Range("E15").Value = "  +11.18%  "
$ws.Range("D16").Value = "3.524.71"
$ws.Range("E16").Value = "  +11.36%  "
$ws.Range("D17").Value = "68.652.80"
$ws.Range("E17").Value = "  +11.01%  "
$ws.Range("E18").Value = "  +8.35%  "
$ws.Range("D19").Value = "18.42"
$ws.Range("E19").Value = "  +9.91%  "
$ws.Range("D20").Value = "11.92"
$ws.Range("E20").Value = "  +12.13%  "
$ws.Range("E21").Value = "  +9.03%  "
$ws.Range("D22").Value = "407.80"
$ws.Range("E22").Value = "  +14.44%  "
$ws.Range("D23").Value = "3.97"
$ws.Range("E23").Value = "  +10.02%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "84.95"
$ws.Range("E24").Value = "  +8.75%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "11.71"
$ws.Range("E25").Value = "  +10.44%  "
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "4.22"
$ws.Range("E26").Value = "  +11.58%  "
$ws.Range("E27").Value = "  +13.20%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").Value = "11.93"
$ws.Range("E29").Value = "  +9.51%  "
$ws.Range("D30").Value = "8.64"
$ws.Range("E30").Value = "  +8.50%  "
$ws.Range("D31").Value = "30.59"
$ws.Range("E31").Value = "  +10.94%  "
$ws.Range("D32").Value = "681.89"
$ws.Range("E32").Value = "  +11.76%  "
$ws.Range("D33").Value = "6.89"
$ws.Range("E33").Value = "  +9.38%  "
$ws.Range("E34").Value = "  +8.30%  "
$ws.Range("E35").Value = "  +10.43%  "
$ws.Range("D36").Value = "60.45"
$ws.Range("E36").Value = "  +6.95%  "
$ws.Range("D37").Value = "0.0₃0855"
$ws.Range("E37").Value = "  +30.31%  "
$ws.Range("E38").Value = "  +9.91%  "
$ws.Range("D39").Value = "0.405"
$ws.Range("E39").Value = "  +9.52%  "
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.08%  "
$ws.Range("D41").Value = "3.41"
$ws.Range("E41").Value = "  +27.48%  "
$ws.Range("D42").Value = "0.134"
$ws.Range("E42").Value = "  +13.12%  "
$ws.Range("E43").Value = "  +17.31%  "
$ws.Range("D44").Value = "3.02"
$ws.Range("E44").Value = "  +18.20%  "
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("D46").Value = "3.047.78"
$ws.Range("E46").Value = "  +9.67%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "9.39"
$ws.Range("E47").Value = "  +28.44%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0423"
$ws.Range("E48").Value = "  +11.84%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.75"
$ws.Range("E49").Value = "  +7.44%  "
$ws.Range("B50").Value = "ApeXProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D50").Value = "3.18"
$ws.Range("E50").Value = "  +11.23%  "
$ws.Range("E51").Value = "  +9.17%  "
